$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (new TPM values)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3518616666666667
$ws.Range("H2").Value = 1.055585
$ws.Range("I2").Value = 0.5958054833396739
$ws.Range("J2").Value = 0.5958054833396738
$ws.Range("M2").Value = 0.08268033333333334
$ws.Range("Q2").Value = 0.02909203988722223
$ws.Range("R2").Value = 0.261828358985
$ws.Range("S2").Value = 0.5958054833396739
$ws.Range("T2").Value = 0.5958054833396738

# Row 3 updates (new TPM values)
$ws.Range("H3").Value = 0.716109
$ws.Range("I3").Value = 0.4041945166603262
$ws.Range("J3").Value = 0.4041945166603262
$ws.Range("M3").Value = 0.08268033333333334
$ws.Range("S3").Value = 0.4041945166603262
$ws.Range("T3").Value = 0.4041945166603262
